$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 395 - this shifts all rows from 395 downward
# one row down (old 395 -> new 396, ..., old 424 -> new 425), matching the
# target diff where the dataset dimension grows from R424 to R425.
$ws.Rows("395:395").Insert()

# Populate the newly inserted row 395 with the new weekly price record.
# All descriptive columns mirror the surrounding rows for this same
# market / product (Agrícola del Norte S.A. de Arica - Brócoli).
$ws.Cells.Item(395, 1).Value = 1
$ws.Cells.Item(395, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(395, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(395, 4).Value = 44783
$ws.Cells.Item(395, 5).Value = 15
$ws.Cells.Item(395, 6).Value = 100112023
$ws.Cells.Item(395, 7).Value = "Brócoli"
$ws.Cells.Item(395, 8).Value = "Sin especificar"
$ws.Cells.Item(395, 9).Value = "Tercera"
$ws.Cells.Item(395, 10).Value = 800
$ws.Cells.Item(395, 11).Value = 1200
$ws.Cells.Item(395, 12).Value = 1300
$ws.Cells.Item(395, 13).Value = 1250
$ws.Cells.Item(395, 14).Value = "$/unidad"
$ws.Cells.Item(395, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(395, 16).Value = 1250
$ws.Cells.Item(395, 17).Value = 1
$ws.Cells.Item(395, 18).Value = "Hortaliza"
